# "Generate Report for Handback"
#
# This localization-status report gets two files (f2cd8a88-....md and
# ffffc3757e19-....md) handed back: zh-cn finished first, de-de a little
# later. For every language sheet the report fills in the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns for
# both data rows (the sheet mirrors row 2's handoff file onto row 3, same
# as the rest of the sheet already does), links the "Latest Target File"
# cell back to the source doc, and flips the Overview/per-language
# "Status" column from "Ready for handoff" to "Handed back: in sync with
# en-US".

$wb = $excel.ActiveWorkbook

$srcDocUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/551e5bfefc155e3cbcc1768dd4abbeb03ea106f6/e2e/f2cd8a88-8020-4b70-ab5b-5a958ea4bffe.md"
$srcDocName = "f2cd8a88-8020-4b70-ab5b-5a958ea4bffe.md"
$newStatus  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns (E, F) and let
# the Status text refresh (same shared cell text as below).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17
$wsOverview.Range("E2:F3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): same shape of edit on each, just a
# different hand-off file / datetime.
# ---------------------------------------------------------------------
function Update-LanguageSheet {
    param($ws, [string]$handoffXlf, [string]$handbackDateTime)

    # Column widths: Status (C) and Latest Target/Handback File (I, J)
    $ws.Range("C1").ColumnWidth = 29.17
    $ws.Range("I1").ColumnWidth = 39.17
    $ws.Range("J1").ColumnWidth = 39.17

    # Status column -> handed back
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File / Latest Handback File / Latest Handback DateTime
    $ws.Range("I2").Value = $srcDocName
    $ws.Range("J2").Value = $handoffXlf
    $ws.Range("K2").Value = $handbackDateTime

    $ws.Range("I3").Value = $srcDocName
    $ws.Range("J3").Value = $handoffXlf
    $ws.Range("K3").Value = $handbackDateTime

    # Re-create the hyperlinks so the new ones land right after row 2's
    # existing hyperlink and before row 3's (matching the handback order
    # the files were produced in): delete + re-add row 3's hyperlink so
    # it gets pushed after the new row-2 link, then add the row-3 link.
    $hlRow3 = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$3') {
            $hlRow3 = $hl
        }
    }
    $row3Address = $hlRow3.Address
    $row3Display = $hlRow3.TextToDisplay
    $hlRow3.Delete()

    $ws.Hyperlinks.Add($ws.Range("I2"), $srcDocUrl, "", "", $srcDocName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $row3Address, "", "", $row3Display)
    $ws.Hyperlinks.Add($ws.Range("I3"), $srcDocUrl, "", "", $srcDocName)

    # Hyperlinks.Add() re-stamps the cell with Excel's generic built-in
    # "Hyperlink" style (theme color); re-apply the workbook's existing
    # custom hyperlink look (blue + underline, same as column A) so I2/I3
    # (and the re-created A3) match the rest of the sheet.
    $ws.Range("I2:I3").Font.Underline = 2
    $ws.Range("I2:I3").Font.Color = 15570276
    $ws.Range("A3").Font.Underline = 2
    $ws.Range("A3").Font.Color = 15570276
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZhCn "f2cd8a88-8020-4b70-ab5b-5a958ea4bffe.6a42664e605b5a09a9b0fb1f6178ead5077864f9.zh-cn.xlf" "2016-09-09 12:48:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDeDe "f2cd8a88-8020-4b70-ab5b-5a958ea4bffe.6a42664e605b5a09a9b0fb1f6178ead5077864f9.de-de.xlf" "2016-09-09 12:49:06"
